# Remove columns that do not affect simulation outcome:
#   H  - "envelope type"
#   I  - "base material"
#   T  - "wall orientations"
#
# Deleting from right-to-left so earlier column letters stay valid
# while later ones are removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("T:T").EntireColumn.Delete()
$ws.Range("H:I").EntireColumn.Delete()

$ws.Range("A2").Select()
